# Add a new "20200916_Morning_Noise" RunParameters sheet for Run30.
#
# The new sheet is created by copying the existing "20200916_RnPoAlphaEffTest"
# sheet (which already has the standard RunParameters layout/values/styles)
# and placing the copy immediately before it, then renaming the copy. This
# reproduces the target workbook structure: a brand-new worksheet inserted at
# position 4 (sheetId 12), with every sheet after it shifting down by one
# position, while all of the existing sheets/data remain untouched.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("20200916_RnPoAlphaEffTest")

# Copy the source sheet to right before itself -- this places the new sheet
# at index 4, right after "20200915_Night_AfterFirstRnInje".
$source.Copy($source)

$newSheet = $wb.Worksheets.Item("20200916_RnPoAlphaEffTest (2)")
$newSheet.Name = "20200916_Morning_Noise"
